$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 17, shifting existing rows 17-19 down to 18-20
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with the new task
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = "Front-end: Data charts"
$ws.Cells.Item(17, 3).Value = "LMS v3"
$ws.Cells.Item(17, 4).Value = 4

# Clear any formatting/content that may have been copied into the new row for columns E:I
$ws.Range("E17:I17").Clear()

# Fix up the S/N column for the rows that shifted down
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(20, 1).Value = 10

# Update the sheet view: remove the frozen/scrolled topLeftCell and change selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E10").Select() | Out-Null
